# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets,
# and append a new event row (row 41) to both sheets.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value (column F = "想去人数")
$fUpdates = @{
    3  = 1327
    5  = 259
    7  = 93
    9  = 177
    11 = 4460
    12 = 6725
    13 = 40
    14 = 56
    16 = 565
    18 = 4102
    19 = 469
    20 = 70
    21 = 50
    22 = 2683
    25 = 163
    26 = 349
    27 = 350
    29 = 217
    30 = 31
    31 = 1614
    32 = 1015
    34 = 124
    35 = 78
    36 = 537
    38 = 11
    39 = 86
    40 = 628
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }

    # Append new row 41 with the new event entry
    $ws.Cells.Item(41, 1).Value = 40
    # Copy A40's formatting (bold/centered/bordered) onto A41, same as every
    # other row's first ("序号") column in this sheet.
    $ws.Cells.Item(40, 1).Copy()
    $ws.Cells.Item(41, 1).PasteSpecial(-4122)
    # Leading apostrophe forces Excel to keep this date-shaped value as text
    # (matches how the other "开始时间" cells in this column are stored).
    $ws.Cells.Item(41, 2).Value = "'2024-09-15"
    $ws.Cells.Item(41, 3).Value = "南昌·第一届哥布林动漫游戏展——开学季&贺中秋"
    $ws.Cells.Item(41, 4).Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    $ws.Cells.Item(41, 5).Value = "2024.09.15 10:00-09.16 18:00"
    $ws.Cells.Item(41, 6).Value = 6
    $ws.Cells.Item(41, 7).Value = 55
    $ws.Cells.Item(41, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89240"
    $ws.Cells.Item(41, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/pixnzm5p1720496832036.jpeg"
}
